$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. trailing zeros, dot-grouped
# thousands) are preserved verbatim instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.980.47"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.751.05"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.15"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.28"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.748.21"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.68"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000247"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.378.98"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.742.56"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.967.28"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.69"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.32"
$ws.Range("E21").Value = "  +5.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.13"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.79"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.29"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.15"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.42"
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.899.25"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.69"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.686.90"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.93"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  +3.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.10"
$ws.Range("E41").Value = "  +8.18%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.325"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.60"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "423.60"
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.44"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.16"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.60"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.784.59"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  +7.07%  "
